# Melhoria do fluxo alternativo 7
# This script rewires the "tipo de cliente / quantidade" steps (step 7 & 8)
# across the five test cases (TC1..TC5) so that each test case now exercises
# a distinct client-type / quantity combination, and trims the now-redundant
# tail (steps 9-12 + postcondition) off TC5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- TC1 (rows 10-22): step 7/8 ---
$ws.Range("B16").Value = "Usuário do Sistema mantém seleção padrão do tipo de cliente A"
$ws.Range("B17").Value = "Usuário do Sistema informa a quantidade de produtos"
$ws.Range("D17").Value = "SYSTEM registra a quantidade informada"

# --- TC2 (rows 29-41): step 7/8 ---
$ws.Range("B35").Value = "Usuário do Sistema altera para tipo de cliente C"
$ws.Range("B36").Value = "Usuário do Sistema informa a quantidade de produtos entre 100 e 999 unidades"
$ws.Range("D36").Value = "SYSTEM aplica fator de desconto para 100 <= quantidade < 1000: Cliente A (0,95), B (0,90), C (0,85)"

# --- TC3 (rows 48-60): step 7/8 ---
$ws.Range("B54").Value = "Usuário do Sistema altera para tipo de cliente B"
$ws.Range("B55").Value = "Usuário do Sistema informa a quantidade de produtos igual ou maior que 1000 unidades"
$ws.Range("D55").Value = "SYSTEM aplica fator de desconto para quantidade >= 1000: Cliente A (1,00), B (0,95), C (0,90)"

# --- TC4 (rows 67-79): step 8 (step 7 unchanged) ---
$ws.Range("B74").Value = "Usuário do Sistema informa a quantidade de produtos menor que 100 unidades"
$ws.Range("D74").Value = "SYSTEM aplica fator de desconto para quantidade < 100: Cliente A (0,90), B (0,85), C (0,80)"

# --- TC5 (rows 86-93): step 7/8 ---
$ws.Range("B92").Value = "Usuário do Sistema altera para tipo de cliente B"
$ws.Range("B93").Value = "Usuário do Sistema informa a quantidade de produtos menor ou igual a zero"
$ws.Range("D93").Value = "SYSTEM exibe mensagem 'A quantidade informada deve ser maior ou igual a 01 (um)!' (MSG002)"

# TC5 now ends right after the negative-quantity validation (row 93); the
# trailing steps 9-12 + postcondition (rows 94-98), which just duplicated
# the other test cases' success flow, are removed entirely.
$ws.Range("A94:F98").EntireRow.Delete()
